# Insert a new data row at row 381 (pushing existing rows 381..475 down to 382..476)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(381).Insert()

$ws.Range("A381").Value2 = 9
$ws.Range("B381").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C381").Value2 = "Metropolitana"
$ws.Range("D381").Value2 = 44785
$ws.Range("D381").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E381").Value2 = 13
$ws.Range("F381").Value2 = 100112012
$ws.Range("G381").Value2 = "Espinaca"
$ws.Range("H381").Value2 = "Sin especificar"
$ws.Range("I381").Value2 = "Primera"
$ws.Range("J381").Value2 = 160
$ws.Range("K381").Value2 = 7000
$ws.Range("L381").Value2 = 8000
$ws.Range("M381").Value2 = 7500
$ws.Range("N381").Value2 = "`$/cuna 10 kilos"
$ws.Range("O381").Value2 = "Provincia de Chacabuco"
$ws.Range("P381").Value2 = 750
$ws.Range("Q381").Value2 = 10
$ws.Range("R381").Value2 = "Hortaliza"
